$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ManageListings")

# Update start date in D2
$ws.Range("D2").Value = 45083

# Update F2/G2 to become text values with the new "am/pm" style formatting
$ws.Range("F2").Value = "12.24.am"
$ws.Range("G2").Value = "19.24.pm"
